$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Test Data")

# A "format donor" cell already carrying the default (unformatted) style
# used throughout the data rows of Sheet1 - used below to restore the
# original style after forcing text storage for numeric-looking values.
$styleDonor = $ws1.Range("A2")

function Set-TextValue($range, [string]$text) {
    # Excel normally auto-converts digit-only strings back into numbers.
    # Force text storage by switching to the Text number format first...
    $range.NumberFormat = "@"
    $range.Value = $text
    # ...then restore the plain/default formatting (style) by copying it
    # from a cell that already has the sheet's normal (unstyled) look,
    # without touching the freshly-typed text value.
    $styleDonor.Copy()
    $range.PasteSpecial(-4122)
}

# Append three new credential rows to Sheet1 (rows 53-55)
Set-TextValue $ws1.Range("A53") "7980000082"
$ws1.Range("B53").Value = "test21726"
$ws1.Range("C53").Value = "test21726@gmail.com"
$ws1.Range("D53").Value = "SoftSuave128316"

Set-TextValue $ws1.Range("A54") "7980000083"
$ws1.Range("B54").Value = "test21726"
$ws1.Range("C54").Value = "test21726@gmail.com"
$ws1.Range("D54").Value = "SoftSuave128316"

Set-TextValue $ws1.Range("A55") "7980000084"
$ws1.Range("B55").Value = "test21726"
$ws1.Range("C55").Value = "test21726@gmail.com"
$ws1.Range("D55").Value = "SoftSuave128316"

# Mark the corresponding mobile numbers as used on the "Test Data" sheet
$ws2.Range("B83").Value = "used"
$ws2.Range("B84").Value = "used"
$ws2.Range("B85").Value = "used"
